$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 with Exp 19 HTM parameters, matching the style/format of
# the preceding rows (copy row 14 formatting down into row 15 first).
$ws.Range("A14:F14").Copy() | Out-Null
$ws.Range("A15:F15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A15").Value = "Exp 19"
$ws.Range("B15").Value = 0.45
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "Local"
$ws.Range("E15").Value = -1
$ws.Range("F15").Value = "Exp 19.png"

# Move the active selection to F16, matching the final cursor position
# recorded after the edit.
$ws.Range("F16").Select() | Out-Null
